$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds header "K" (Strike count). Values regenerated to new K values.
$newValues = @{
    2  = 2
    3  = 1
    4  = 0
    5  = 2
    6  = 3
    7  = 4
    8  = 7
    9  = 6
    10 = 3
    11 = 3
    12 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
